$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.056.23"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "3.061.12"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.058.37"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("E12").Value = "  -3.36%  "
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "3.571.97"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "66.032.90"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("E18").Value = "  -3.35%  "
$ws.Range("D19").Value = "3.071.13"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "481.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.08%  "
$ws.Range("E26").Value = "  -3.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("D34").Value = "0.0₃0892"
$ws.Range("E34").Value = "  -6.32%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "47.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.938"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.51%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -5.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.297"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.59%  "
$ws.Range("D43").Value = "2.759.55"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "134.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0340"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "359.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("E51").Value = "  -3.20%  "
